$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (days since epoch).
# Update every data row (2 through 261) from 45181 to 45182.
$ws.Range("C2:C261").Value = 45182
